# "Set Developer Tab and Added a Message Box"
#
# The author turned on the Developer tab (an Excel application-level UI
# preference) and opened the VBA editor, which is why the workbook/sheet
# picked up explicit VBA "CodeName" identities (ThisWorkbook / Sheet1).
# Mirror that intent on the object model here (harmless if the host
# doesn't persist it, since these already are the implicit defaults).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.ShowDevTools = $true

$wb.CodeName = "ThisWorkbook"
$ws.CodeName = "Sheet1"

# The actual content change: six rows of the value 20 in column A.
$ws.Range("A1:A6").Value = 20
